$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows: B2 and B3 changed from 1 to 0
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0

# Add new rows 5-7
$ws.Range("A5").Value = "ISTAT_12_323_DF_DCCV_IMPDEP_1"
$ws.Range("B5").Value = 1

$ws.Range("A6").Value = "ISTAT_12_323_DF_DCCV_IMPDEP_2"
$ws.Range("B6").Value = 1

$ws.Range("A7").Value = "ISTAT_12_60_DF_DCCV_CONSACQUA_1"
$ws.Range("B7").Value = 1

# Widen column A to fit new longer text
$ws.Columns.Item(1).ColumnWidth = 43.5

# Set the selection to C13 to match the saved view state
$ws.Range("C13").Select()
